$wb = $excel.ActiveWorkbook

# Rename the first sheet (tab name) to standardize naming
$ws = $wb.Worksheets.Item("total_reg_and_ballots")
$ws.Name = "total_reg_and_cast"

# Update the selection on that sheet from B38 to A47
$ws.Activate()
$ws.Range("A47").Select()
